$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and volume change (E) columns with latest scraped values.
# Each cell is forced to Text format before assignment so values such as
# '24.828.42' or '1.001' are not auto-converted to numbers/dates by Excel,
# then the style is reset back to Normal so no stray style index is left on the cell.
$updates = @{
    'D2' = '24.828.42'
    'E2' = '  +1.72%  '
    'D3' = '1.708.58'
    'E3' = '  +1.82%  '
    'E4' = '  -0.05%  '
    'D5' = '311.68'
    'E5' = '  +1.79%  '
    'D6' = '0.9986'
    'E6' = '  +0.07%  '
    'D7' = '0.3763'
    'E7' = '  +1.35%  '
    'D8' = '49.71'
    'E8' = '  +2.93%  '
    'D9' = '0.3451'
    'E9' = '  +0.28%  '
    'E10' = '  +2.54%  '
    'D11' = '0.07555'
    'E11' = '  +3.97%  '
    'D12' = '1.001'
    'E12' = '  -0.09%  '
    'D13' = '21.18'
    'E13' = '  +3.86%  '
    'D14' = '6.314'
    'E14' = '  +3.01%  '
    'D15' = '7.078'
    'E15' = '  +4.98%  '
    'D16' = '1.712.53'
    'E16' = '  +2.07%  '
    'D17' = '0.00001135'
    'E17' = '  +2.57%  '
    'D18' = '0.06722'
    'E18' = '  +0.09%  '
    'D19' = '0.9990'
    'E19' = '  +0.08%  '
    'D20' = '85.04'
    'E20' = '  +4.66%  '
    'D21' = '17.37'
    'E21' = '  +5.66%  '
    'D22' = '6.408'
    'E22' = '  +5.04%  '
    'D23' = '13.31'
    'E23' = '  +11.24%  '
    'D24' = '24.840.11'
    'E24' = '  +1.81%  '
    'D25' = '2.458'
    'E25' = '  +0.92%  '
    'D26' = '2.806'
    'E26' = '  +5.25%  '
    'E27' = '  +5.02%  '
    'D28' = '152.03'
    'E28' = '  -0.28%  '
    'D29' = '133.06'
    'E29' = '  +4.65%  '
    'D30' = '1.900.89'
    'E30' = '  +2.13%  '
    'D31' = '1.255'
    'E31' = '  +29.52%  '
    'D32' = '6.952'
    'E32' = '  +9.49%  '
    'D33' = '4.216'
    'E33' = '  +4.74%  '
    'D34' = '13.96'
    'E34' = '  +12.27%  '
    'D35' = '1.796'
    'E35' = '  +6.80%  '
    'D36' = '0.08843'
    'E36' = '  +4.46%  '
    'D37' = '9.410'
    'E37' = '  +5.22%  '
    'D38' = '5.632'
    'E38' = '  +5.37%  '
    'D39' = '0.06696'
    'E39' = '  +3.02%  '
    'D40' = '0.02417'
    'E40' = '  +3.40%  '
    'D41' = '0.2242'
    'E41' = '  +6.17%  '
    'D42' = '1.284'
    'E42' = '  +1.61%  '
    'D43' = '0.6459'
    'E43' = '  +4.41%  '
    'D44' = '0.9988'
    'E44' = '  +0.20%  '
    'D45' = '13.95'
    'E45' = '  +7.14%  '
    'D46' = '0.6178'
    'E46' = '  +3.76%  '
    'D47' = '3.823'
    'E47' = '  +1.45%  '
    'D48' = '2.139'
    'E48' = '  +5.35%  '
    'D49' = '130.40'
    'E49' = '  +2.70%  '
    'D50' = '0.07323'
    'E50' = '  +1.44%  '
    'D51' = '80.26'
    'E51' = '  +6.01%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
